$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "完成情况" (completion status) column for the last group's
# task rows (人员 = 钟崇尧 / 黄宝怡 / 魏仲凯 / 薛洁鹏 / 王晓宇), which were
# previously left blank.
$ws.Range("C51").Value = "未完成、完成一半"
$ws.Range("C52").Value = "已完成"
$ws.Range("C53").Value = "已完成"
$ws.Range("C54").Value = "已完成"
$ws.Range("C55").Value = "已完成"

# Leave the selection where the edits were made, matching the saved
# workbook's cursor position.
[void]$ws.Range("C51").Select()
